# Commit 5 : amelioration import transactions
#
# The transaction log is extended from 3 rows to 9 rows: a new "Achat"
# entry is inserted at the top (pushing the existing rows down), and
# several new "Achat"/"Vente" rows are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, extend the sheet with extra rows (6 new rows) that have the
# same formatting (number format / borders) as the existing data rows,
# by duplicating row 2 into rows 5:10. Values will be overwritten next.
$ws.Range("A2:G2").Copy($ws.Range("A5:G10"))

# Row 2 : new first transaction (Achat, 17 Feb 2026)
$ws.Range("A2").Value = 46070
$ws.Range("B2").Value = "Achat"
$ws.Range("C2").Value = "LU0908500753"
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "EUR"
$ws.Range("G2").Value = "Portfolio Fictif Test"

# Row 3 : former row 2 data (Achat, 10 Feb 2026)
$ws.Range("A3").Value = 46063
$ws.Range("B3").Value = "Achat"
$ws.Range("C3").Value = "LU0908500753"
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = "EUR"
$ws.Range("G3").Value = "Portfolio Fictif Test"

# Row 4 : former row 3 data (Achat, 13 Feb 2026)
$ws.Range("A4").Value = 46066
$ws.Range("B4").Value = "Achat"
$ws.Range("C4").Value = "LU0908500753"
$ws.Range("D4").Value = 50
$ws.Range("E4").Value = 12
$ws.Range("F4").Value = "EUR"
$ws.Range("G4").Value = "Portfolio Fictif Test"

# Row 5 : new row (Achat, 15 Feb 2026)
$ws.Range("A5").Value = 46068
$ws.Range("B5").Value = "Achat"
$ws.Range("C5").Value = "LU0908500753"
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = 20
$ws.Range("F5").Value = "EUR"
$ws.Range("G5").Value = "Portfolio Fictif Test"

# Row 6 : new row (Vente, 20 Feb 2026)
$ws.Range("A6").Value = 46073
$ws.Range("B6").Value = "Vente "
$ws.Range("C6").Value = "LU0908500753"
$ws.Range("D6").Value = 40
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = "EUR"
$ws.Range("G6").Value = "Portfolio Fictif Test"

# Row 7 : former row 4 data (Vente, 16 Feb 2026)
$ws.Range("A7").Value = 46069
$ws.Range("B7").Value = "Vente "
$ws.Range("C7").Value = "LU0908500753"
$ws.Range("D7").Value = 40
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = "EUR"
$ws.Range("G7").Value = "Portfolio Fictif Test"

# Row 8 : new row (Achat, 19 Feb 2026)
$ws.Range("A8").Value = 46072
$ws.Range("B8").Value = "Achat"
$ws.Range("C8").Value = "LU0908500753"
$ws.Range("D8").Value = 50
$ws.Range("E8").Value = 12
$ws.Range("F8").Value = "EUR"
$ws.Range("G8").Value = "Portfolio Fictif Test"

# Row 9 : new row (Achat, 22 Feb 2026)
$ws.Range("A9").Value = 46075
$ws.Range("B9").Value = "Achat"
$ws.Range("C9").Value = "LU0908500753"
$ws.Range("D9").Value = 50
$ws.Range("E9").Value = 12
$ws.Range("F9").Value = "EUR"
$ws.Range("G9").Value = "Portfolio Fictif Test"

# Row 10 : new row (Achat, 26 Feb 2026)
$ws.Range("A10").Value = 46079
$ws.Range("B10").Value = "Achat"
$ws.Range("C10").Value = "LU0908500753"
$ws.Range("D10").Value = 50
$ws.Range("E10").Value = 20
$ws.Range("F10").Value = "EUR"
$ws.Range("G10").Value = "Portfolio Fictif Test"

# Reflect the author's final on-screen selection (whole row 9 selected).
$ws.Rows(9).Select()
